# Auto update Excel log: append newly logged sensor readings (2026-02-06)
# to the PIR, Humidity and Temperature sheets.
$wb = $excel.ActiveWorkbook

# --- PIR sheet: append rows 513-525 ---
$ws = $wb.Worksheets.Item("PIR")
$ws.Range("A513:A525").NumberFormat = "@"  # keep as text, not auto-converted
$newRows = @(
    @('2026-02-06', '10:22:35', '10:00', 'Bathroom', 'No Motion', 'Inactive'),
    @('2026-02-06', '10:22:38', '10:00', 'Bathroom', 'No Motion', 'Inactive'),
    @('2026-02-06', '10:22:39', '10:00', 'Bathroom', 'Motion Detected', 'Active'),
    @('2026-02-06', '10:22:45', '10:00', 'Bathroom', 'No Motion', 'Inactive'),
    @('2026-02-06', '10:22:50', '10:00', 'Bathroom', 'No Motion', 'Inactive'),
    @('2026-02-06', '10:22:55', '10:00', 'Bathroom', 'No Motion', 'Inactive'),
    @('2026-02-06', '10:23:00', '10:00', 'Bathroom', 'No Motion', 'Inactive'),
    @('2026-02-06', '10:23:05', '10:00', 'Bathroom', 'No Motion', 'Inactive'),
    @('2026-02-06', '10:23:10', '10:00', 'Bathroom', 'No Motion', 'Inactive'),
    @('2026-02-06', '10:23:15', '10:00', 'Bathroom', 'No Motion', 'Inactive'),
    @('2026-02-06', '10:23:20', '10:00', 'Bathroom', 'No Motion', 'Inactive'),
    @('2026-02-06', '10:23:25', '10:00', 'Bathroom', 'No Motion', 'Inactive'),
    @('2026-02-06', '10:23:30', '10:00', 'Bathroom', 'No Motion', 'Inactive')
)
$r = 513
foreach ($row in $newRows) {
    for ($c = 1; $c -le 6; $c++) {
        $ws.Cells.Item($r, $c).Value = $row[$c - 1]
    }
    $r++
}

# --- Humidity sheet: append rows 360-368 ---
$ws = $wb.Worksheets.Item("Humidity")
$ws.Range("A360:A368").NumberFormat = "@"  # keep as text, not auto-converted
$ws.Range("E360:E368").NumberFormat = "@"  # keep as text, not auto-converted
$newRows = @(
    @('2026-02-06', '10:22:36', '10:00', 'Bathroom', '68.3%', 'Active'),
    @('2026-02-06', '10:22:46', '10:00', 'Bathroom', '68.3%', 'Active'),
    @('2026-02-06', '10:22:51', '10:00', 'Bathroom', '68.3%', 'Active'),
    @('2026-02-06', '10:22:56', '10:00', 'Bathroom', '68.4%', 'Active'),
    @('2026-02-06', '10:23:01', '10:00', 'Bathroom', '68.3%', 'Active'),
    @('2026-02-06', '10:23:06', '10:00', 'Bathroom', '68.2%', 'Active'),
    @('2026-02-06', '10:23:17', '10:00', 'Bathroom', '68.1%', 'Active'),
    @('2026-02-06', '10:23:22', '10:00', 'Bathroom', '68.1%', 'Active'),
    @('2026-02-06', '10:23:27', '10:00', 'Bathroom', '68.1%', 'Active')
)
$r = 360
foreach ($row in $newRows) {
    for ($c = 1; $c -le 6; $c++) {
        $ws.Cells.Item($r, $c).Value = $row[$c - 1]
    }
    $r++
}

# --- Temperature sheet: append rows 360-368 ---
$ws = $wb.Worksheets.Item("Temperature")
$ws.Range("A360:A368").NumberFormat = "@"  # keep as text, not auto-converted
$newRows = @(
    @('2026-02-06', '10:22:37', '10:00', 'Bathroom', '28.2C', 'Active'),
    @('2026-02-06', '10:22:48', '10:00', 'Bathroom', '28.2C', 'Active'),
    @('2026-02-06', '10:22:53', '10:00', 'Bathroom', '28.2C', 'Active'),
    @('2026-02-06', '10:22:58', '10:00', 'Bathroom', '28.3C', 'Active'),
    @('2026-02-06', '10:23:03', '10:00', 'Bathroom', '28.2C', 'Active'),
    @('2026-02-06', '10:23:08', '10:00', 'Bathroom', '28.3C', 'Active'),
    @('2026-02-06', '10:23:18', '10:00', 'Bathroom', '28.3C', 'Active'),
    @('2026-02-06', '10:23:23', '10:00', 'Bathroom', '28.3C', 'Active'),
    @('2026-02-06', '10:23:28', '10:00', 'Bathroom', '28.3C', 'Active')
)
$r = 360
foreach ($row in $newRows) {
    for ($c = 1; $c -le 6; $c++) {
        $ws.Cells.Item($r, $c).Value = $row[$c - 1]
    }
    $r++
}
